$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove Sheet2 entirely (results for a removed run were dropped)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Delete()

$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 4 (question 3): CWM column corrected from 1 to 0
$ws1.Range("F4").Value = 0

# Row 10 (question 9): CWM column corrected from 1 to 0
$ws1.Range("F10").Value = 0

# Row 11 (question 10): CWM column corrected from 1 to 0
$ws1.Range("F11").Value = 0

# Row 14: drop the stray annotation cell in M14 (error-type note tied to
# the now-removed Sheet2 breakdown)
$ws1.Range("M14").Clear()

# Row 45 (question 44): CWM column corrected from 1 to 0
$ws1.Range("F45").Value = 0

# Row 53 (question 52): replace the stray text annotation with the
# numeric result value
$ws1.Range("F53").ClearFormats()
$ws1.Range("F53").Value = 0

# Row 55 (question 54): CWM column corrected from 1 to 0
$ws1.Range("F55").Value = 0

# Restore the view to an unscrolled state with the new selection
$ws1.Range("E56").Select()
